$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin name / link columns) -- never numeric-looking, safe to set directly.
$textUpdates = @(
    @("B38", 'InjectiveProtocol'),
    @("C38", 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @("B39", 'Dai'),
    @("C39", 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
    @("B40", 'PEPE'),
    @("C40", 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'),
    @("B42", 'Kaspa'),
    @("C42", 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'),
    @("B43", 'dogwifhat'),
    @("C43", 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif')
)
foreach ($pair in $textUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Price / Volume columns hold numeric-looking strings (e.g. "0.623", "1.00") that Excel
# would otherwise auto-convert to numbers on assignment. Force Text format, assign, then
# restore the Normal style so no residual number-format diff is left behind.
$forcedTextUpdates = @(
    @("D2", '70.842.59'),
    @("E2", '  +0.98%  '),
    @("D3", '3.589.14'),
    @("E3", '  +0.35%  '),
    @("E4", '  +0.01%  '),
    @("D5", '587.73'),
    @("E5", '  +1.46%  '),
    @("D6", '186.80'),
    @("E6", '  +0.15%  '),
    @("D7", '3.578.79'),
    @("E7", '  +0.22%  '),
    @("D8", '0.623'),
    @("E8", '  +0.60%  '),
    @("E9", '  +0.00%  '),
    @("D10", '0.214'),
    @("E10", '  +17.22%  '),
    @("D11", '0.653'),
    @("E11", '  +0.33%  '),
    @("D12", '54.49'),
    @("E12", '  -1.18%  '),
    @("D13", '0.0000320'),
    @("E13", '  +5.10%  '),
    @("D14", '9.58'),
    @("E14", '  +0.41%  '),
    @("D15", '4.158.22'),
    @("E15", '  +0.14%  '),
    @("D16", '19.63'),
    @("E16", '  -0.15%  '),
    @("D17", '70.829.60'),
    @("E17", '  +1.10%  '),
    @("D18", '3.589.85'),
    @("E18", '  +0.21%  '),
    @("D19", '578.48'),
    @("E19", '  +17.43%  '),
    @("D20", '12.44'),
    @("E20", '  -0.99%  '),
    @("E21", '  +0.00%  '),
    @("E22", '  -1.71%  '),
    @("D23", '17.86'),
    @("E23", '  -8.04%  '),
    @("D24", '4.67'),
    @("E24", '  +6.19%  '),
    @("D25", '4.94'),
    @("E25", '  -0.61%  '),
    @("D26", '95.63'),
    @("E26", '  -1.15%  '),
    @("D27", '11.49'),
    @("E27", '  +0.16%  '),
    @("D28", '2.96'),
    @("E28", '  -0.02%  '),
    @("D29", '9.17'),
    @("E29", '  -1.86%  '),
    @("D30", '32.36'),
    @("E30", '  +2.07%  '),
    @("D31", '7.35'),
    @("E31", '  -5.74%  '),
    @("D32", '12.44'),
    @("E32", '  +2.72%  '),
    @("D33", '65.10'),
    @("E33", '  -1.22%  '),
    @("E34", '  -0.61%  '),
    @("D35", '3.38'),
    @("E35", '  +5.08%  '),
    @("D36", '564.36'),
    @("E36", '  -2.79%  '),
    @("D37", '0.419'),
    @("E37", '  +1.74%  '),
    @("D38", '37.82'),
    @("E38", '  -2.97%  '),
    @("D39", '1.00'),
    @("E39", '  +0.12%  '),
    @("D40", '0.0₃0798'),
    @("E40", '  +0.46%  '),
    @("D41", '3.397.98'),
    @("E41", '  +6.70%  '),
    @("D42", '0.136'),
    @("E42", '  +0.29%  '),
    @("D43", '3.14'),
    @("E43", '  -0.97%  '),
    @("D44", '3.40'),
    @("E44", '  -2.34%  '),
    @("D45", '3.59'),
    @("E45", '  -0.16%  '),
    @("D46", '0.0449'),
    @("E46", '  +1.51%  '),
    @("E47", '  -2.95%  '),
    @("D48", '9.37'),
    @("E48", '  -1.35%  '),
    @("D49", '0.138'),
    @("E49", '  +1.25%  '),
    @("D50", '0.998'),
    @("E50", '  -0.20%  '),
    @("D51", '1.43'),
    @("E51", '  -10.62%  ')
)
foreach ($pair in $forcedTextUpdates) {
    $ref = $pair[0]
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $pair[1]
    $ws.Range($ref).Style = "Normal"
}
